$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching style of existing header row (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I ("I0") and J ("IF"), rows 2-59
$iValues = @(8,8,8,3,8,6,8,6,8,11,8,8,8,8,1,4,7,3,8,5,7,9,7,1,6,7,9,7,8,7,7,6,4,9,7,5,9,1,6,6,6,5,4,6,5,6,6,9,8,7,7,6,8,7,7,5,8,4)
$jValues = @(8,8,9,5,8,8,8,6,8,11,8,8,8,9,3,4,7,4,8,5,8,9,7,2,6,9,9,7,8,7,7,7,4,9,8,6,9,3,6,6,7,5,4,7,5,6,7,9,8,7,8,6,8,7,7,5,9,4)

$startRow = 2
for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $startRow + $idx
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
